$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# Remove the hyperlink on I11 (cost_per_mj row) before deleting the row.
$ws.Range("I11").Hyperlinks.Delete()

# Delete row 14 (co2_conversion_per_mj) first so the row-11 index below stays valid.
$ws.Rows.Item(14).Delete()
# Delete row 11 (cost_per_mj) - the rest of the rows shift up to fill the gap.
$ws.Rows.Item(11).Delete()

Write-Host "done"
